$d = $word.ActiveDocument

# First paragraph in the document holds the hidden "**ID__...__ID**" marker
# run followed by a trailing space-only run.
$p = $d.Paragraphs(1)

# Collapse the marker text + trailing space run into the new marker text
# (no trailing space) in a single run.
$markerRange = $d.Range(0, 32)
$markerRange.Text = "**ID__AFFARS_5343_204_70_3__ID**"

# Indent left from 120 (6pt) to 225 twips (11.25pt).
$p.Format.LeftIndent = 11.25

# Add a paragraph border with 5-twip spacing on all sides (matching the
# border already used later in the document).
$borders = $p.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
